# Reading Modality, StudyDescription and StudyDate from DICOM files.
# Populate the newly-added H (Modality), I (StudyDescription) and J (StudyDate)
# columns for each scan data row in the "Files" sheet, using the StudyDate
# that was already present in column F for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> StudyDate (same value already present in column F for that row)
$rows = @{
    12 = "20200312"
    14 = "20200312"
    17 = "20200303"
    21 = "20190115"
    24 = "20200623"
    26 = "20210414"
}

foreach ($r in $rows.Keys) {
    $ws.Range("H$r").Value = "CT"
    $ws.Range("I$r").Value = "CT1 abdomen"

    # Write the StudyDate as text (matching the existing F-column shared
    # string) rather than letting it be auto-detected as a number: build it
    # via a text formula, then convert the formula to its static value.
    $dateCell = $ws.Range("J$r")
    $dateCell.Formula = "=""" + $rows[$r] + """"
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = $false

# Scroll the view down and move the selection, matching the saved view state.
try { $excel.ActiveWindow.ScrollRow = 13 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
$ws.Range("J23").Select()
